# Updated cryptos list (price + 1h volume change refresh)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $value) {
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = "Normal"
}

$ws.Cells.Item(2, 4).Value = "22.465.19"
$ws.Cells.Item(2, 5).Value = "  +0.16%  "
$ws.Cells.Item(3, 4).Value = "1.572.21"
$ws.Cells.Item(3, 5).Value = "  +0.57%  "
Set-TextValue $ws.Cells.Item(4, 4) "1.001"
$ws.Cells.Item(4, 5).Value = "  -0.16%  "
$ws.Cells.Item(5, 5).Value = "  -0.05%  "
Set-TextValue $ws.Cells.Item(6, 4) "288.65"
$ws.Cells.Item(6, 5).Value = "  +0.07%  "
Set-TextValue $ws.Cells.Item(7, 4) "0.3709"
$ws.Cells.Item(7, 5).Value = "  +1.10%  "
Set-TextValue $ws.Cells.Item(8, 4) "48.36"
$ws.Cells.Item(8, 5).Value = "  -2.94%  "
Set-TextValue $ws.Cells.Item(9, 4) "0.3313"
$ws.Cells.Item(9, 5).Value = "  -1.39%  "
Set-TextValue $ws.Cells.Item(10, 4) "1.135"
$ws.Cells.Item(10, 5).Value = "  +0.27%  "
Set-TextValue $ws.Cells.Item(11, 4) "0.07501"
$ws.Cells.Item(11, 5).Value = "  +0.51%  "
Set-TextValue $ws.Cells.Item(12, 4) "1.001"
$ws.Cells.Item(12, 5).Value = "  -0.15%  "
Set-TextValue $ws.Cells.Item(13, 4) "20.75"
$ws.Cells.Item(13, 5).Value = "  -0.77%  "
Set-TextValue $ws.Cells.Item(14, 4) "5.937"
$ws.Cells.Item(14, 5).Value = "  -0.48%  "
Set-TextValue $ws.Cells.Item(15, 4) "6.869"
$ws.Cells.Item(15, 5).Value = "  -0.93%  "
$ws.Cells.Item(16, 4).Value = "1.568.60"
$ws.Cells.Item(16, 5).Value = "  +0.28%  "
Set-TextValue $ws.Cells.Item(17, 4) "0.00001120"
$ws.Cells.Item(17, 5).Value = "  +1.16%  "

# Row 18 becomes TRON (was Litecoin)
$ws.Cells.Item(18, 2).Value = "TRON"
$ws.Cells.Item(18, 3).Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
Set-TextValue $ws.Cells.Item(18, 4) "0.06744"
$ws.Cells.Item(18, 5).Value = "  +0.10%  "

# Row 19 becomes Litecoin (was TRON)
$ws.Cells.Item(19, 2).Value = "Litecoin"
$ws.Cells.Item(19, 3).Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
Set-TextValue $ws.Cells.Item(19, 4) "87.66"
$ws.Cells.Item(19, 5).Value = "  -2.59%  "

$ws.Cells.Item(20, 5).Value = "  -0.05%  "
Set-TextValue $ws.Cells.Item(21, 4) "6.359"
$ws.Cells.Item(21, 5).Value = "  +0.24%  "
Set-TextValue $ws.Cells.Item(22, 4) "16.55"
$ws.Cells.Item(22, 5).Value = "  +2.59%  "
$ws.Cells.Item(23, 5).Value = "  +0.34%  "
$ws.Cells.Item(24, 4).Value = "22.466.14"
$ws.Cells.Item(24, 5).Value = "  +0.19%  "
Set-TextValue $ws.Cells.Item(25, 4) "2.394"
$ws.Cells.Item(25, 5).Value = "  -0.01%  "
Set-TextValue $ws.Cells.Item(26, 4) "2.582"
$ws.Cells.Item(26, 5).Value = "  -1.14%  "
Set-TextValue $ws.Cells.Item(27, 4) "153.67"
$ws.Cells.Item(27, 5).Value = "  +3.06%  "
$ws.Cells.Item(28, 5).Value = "  -0.05%  "
Set-TextValue $ws.Cells.Item(29, 4) "5.019"
$ws.Cells.Item(29, 5).Value = "  -0.55%  "
Set-TextValue $ws.Cells.Item(30, 4) "124.61"
$ws.Cells.Item(30, 5).Value = "  +0.75%  "
$ws.Cells.Item(31, 4).Value = "1.747.47"
$ws.Cells.Item(31, 5).Value = "  +0.54%  "
Set-TextValue $ws.Cells.Item(32, 4) "1.063"
$ws.Cells.Item(32, 5).Value = "  +1.19%  "
Set-TextValue $ws.Cells.Item(33, 4) "2.013"
$ws.Cells.Item(33, 5).Value = "  -0.32%  "
Set-TextValue $ws.Cells.Item(34, 4) "6.124"
$ws.Cells.Item(34, 5).Value = "  +0.13%  "
Set-TextValue $ws.Cells.Item(35, 4) "9.790"
$ws.Cells.Item(35, 5).Value = "  +1.94%  "
Set-TextValue $ws.Cells.Item(36, 4) "0.08363"
$ws.Cells.Item(36, 5).Value = "  +1.08%  "
Set-TextValue $ws.Cells.Item(37, 4) "0.02472"
$ws.Cells.Item(37, 5).Value = "  +1.81%  "
$ws.Cells.Item(38, 5).Value = "  +0.48%  "
Set-TextValue $ws.Cells.Item(39, 4) "0.06417"
$ws.Cells.Item(39, 5).Value = "  +0.18%  "
$ws.Cells.Item(40, 5).Value = "  -3.17%  "
Set-TextValue $ws.Cells.Item(41, 4) "5.348"
$ws.Cells.Item(41, 5).Value = "  +0.63%  "
Set-TextValue $ws.Cells.Item(42, 4) "0.6328"
$ws.Cells.Item(42, 5).Value = "  +3.17%  "
Set-TextValue $ws.Cells.Item(43, 4) "11.31"
$ws.Cells.Item(43, 5).Value = "  +1.93%  "
Set-TextValue $ws.Cells.Item(44, 4) "13.79"
$ws.Cells.Item(44, 5).Value = "  -0.21%  "
Set-TextValue $ws.Cells.Item(45, 4) "0.6173"
$ws.Cells.Item(45, 5).Value = "  +7.19%  "
Set-TextValue $ws.Cells.Item(46, 4) "3.771"
$ws.Cells.Item(46, 5).Value = "  +0.25%  "
Set-TextValue $ws.Cells.Item(47, 4) "2.063"
$ws.Cells.Item(47, 5).Value = "  +1.74%  "
Set-TextValue $ws.Cells.Item(48, 4) "126.14"
$ws.Cells.Item(48, 5).Value = "  +0.51%  "
Set-TextValue $ws.Cells.Item(49, 4) "1.214"
$ws.Cells.Item(49, 5).Value = "  -0.27%  "
Set-TextValue $ws.Cells.Item(50, 4) "0.07220"
$ws.Cells.Item(50, 5).Value = "  -1.40%  "
Set-TextValue $ws.Cells.Item(51, 4) "76.99"
$ws.Cells.Item(51, 5).Value = "  +2.76%  "
